{"js": "// Office.js (Word JavaScript API) port of the authoring diff:\n//   - the leading date paragraph changes from \"2025-03-28 Friday\" to \"2025-03-29 Saturday\"\n//   - every cell of the 20x5 arithmetic-practice table gets a new problem\n// Only the text inside each run changes; every run/paragraph keeps its\n// original formatting (rFonts/sz/jc, etc.) because we replace text through\n// each paragraph's Range rather than rebuilding the paragraph.\n\nconst NEW_DATE = \"2025-03-29 Saturday\";\nconst NEW_VALUES = [[\"45-8=\", \"16+8=\", \"69+19=\", \"29+39=\", \"19+28=\"], [\"6+65=\", \"36+39=\", \"39+43=\", \"36+55=\", \"90-61=\"], [\"87-58=\", \"16+6=\", \"84+9=\", \"77+7=\", \"20-16=\"], [\"63-55=\", \"13+78=\", \"59+6=\", \"8+74=\", \"77+5=\"], [\"37+59=\", \"49+46=\", \"28+43=\", \"27+7=\", \"19+19=\"], [\"90-51=\", \"60-46=\", \"48-29=\", \"89+9=\", \"67-48=\"], [\"30-19=\", \"91-57=\", \"46+15=\", \"81-47=\", \"52-47=\"], [\"74-9=\", \"17+47=\", \"56-8=\", \"6+58=\", \"62+19=\"], [\"64-8=\", \"15+56=\", \"93-37=\", \"60-2=\", \"83-78=\"], [\"9+52=\", \"47+17=\", \"31-5=\", \"65-56=\", \"26+49=\"], [\"18+77=\", \"43+28=\", \"51-8=\", \"83-18=\", \"62-48=\"], [\"53-39=\", \"58-29=\", \"9+56=\", \"75-69=\", \"62-27=\"], [\"85-6=\", \"36-28=\", \"91-58=\", \"91-7=\", \"19+37=\"], [\"18+69=\", \"39+25=\", \"70-62=\", \"26+25=\", \"8+37=\"], [\"97-38=\", \"53-48=\", \"74-47=\", \"27+64=\", \"9+27=\"], [\"14+8=\", \"53-24=\", \"48+25=\", \"38+47=\", \"94-55=\"], [\"69+4=\", \"83-65=\", \"83-58=\", \"26+57=\", \"63-15=\"], [\"49+34=\", \"39+47=\", \"94-25=\", \"26+69=\", \"50-24=\"], [\"18+17=\", \"19+26=\", \"83-79=\", \"16+19=\", \"92-36=\"], [\"42-38=\", \"40-24=\", \"64-36=\", \"95-8=\", \"51-48=\"]];\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph (first paragraph in the document body) ---\nconst dateParagraph = body.paragraphs.getFirst();\ndateParagraph.getRange().insertText(NEW_DATE, Word.InsertLocation.replace);\n\n// --- 2. Update every cell of the (single) table with the new arithmetic values ---\nconst table = body.tables.getFirstOrNullObject();\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body but none was found.\");\n}\n\nfor (let r = 0; r < NEW_VALUES.length; r++) {\n  const row = NEW_VALUES[r];\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCell(r, c);\n    const cellParagraph = cell.body.paragraphs.getFirst();\n    cellParagraph.getRange().insertText(row[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and every math-problem cell in the table,\n# matching the authoring diff (text-only changes; formatting untouched).\n$NewDate = \"2025-03-29 Saturday\"\n\n$NewValues = @(\n    @(\"45-8=\", \"16+8=\", \"69+19=\", \"29+39=\", \"19+28=\"),\n    @(\"6+65=\", \"36+39=\", \"39+43=\", \"36+55=\", \"90-61=\"),\n    @(\"87-58=\", \"16+6=\", \"84+9=\", \"77+7=\", \"20-16=\"),\n    @(\"63-55=\", \"13+78=\", \"59+6=\", \"8+74=\", \"77+5=\"),\n    @(\"37+59=\", \"49+46=\", \"28+43=\", \"27+7=\", \"19+19=\"),\n    @(\"90-51=\", \"60-46=\", \"48-29=\", \"89+9=\", \"67-48=\"),\n    @(\"30-19=\", \"91-57=\", \"46+15=\", \"81-47=\", \"52-47=\"),\n    @(\"74-9=\", \"17+47=\", \"56-8=\", \"6+58=\", \"62+19=\"),\n    @(\"64-8=\", \"15+56=\", \"93-37=\", \"60-2=\", \"83-78=\"),\n    @(\"9+52=\", \"47+17=\", \"31-5=\", \"65-56=\", \"26+49=\"),\n    @(\"18+77=\", \"43+28=\", \"51-8=\", \"83-18=\", \"62-48=\"),\n    @(\"53-39=\", \"58-29=\", \"9+56=\", \"75-69=\", \"62-27=\"),\n    @(\"85-6=\", \"36-28=\", \"91-58=\", \"91-7=\", \"19+37=\"),\n    @(\"18+69=\", \"39+25=\", \"70-62=\", \"26+25=\", \"8+37=\"),\n    @(\"97-38=\", \"53-48=\", \"74-47=\", \"27+64=\", \"9+27=\"),\n    @(\"14+8=\", \"53-24=\", \"48+25=\", \"38+47=\", \"94-55=\"),\n    @(\"69+4=\", \"83-65=\", \"83-58=\", \"26+57=\", \"63-15=\"),\n    @(\"49+34=\", \"39+47=\", \"94-25=\", \"26+69=\", \"50-24=\"),\n    @(\"18+17=\", \"19+26=\", \"83-79=\", \"16+19=\", \"92-36=\"),\n    @(\"42-38=\", \"40-24=\", \"64-36=\", \"95-8=\", \"51-48=\"),\n)\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph (first paragraph in the document) ---\n$d.Paragraphs.Item(1).Range.Text = $NewDate\n\n# --- 2. Update every cell of the (single) table with the new arithmetic values ---\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $NewValues.Count; $r++) {\n    $row = $NewValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n\n"}
